# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-16 13:20:24
#
# Applies the Y2 B25/26 GIT & Liver session-analysis updates:
#  - widen the "Students" column (I)
#  - re-order the "Recorded By" email lists on rows 2, 3 and 15
#  - refresh the Missing/Pending session counters (L7, L8, P15, Q15)
#  - flip the still-not-recorded PHYSIOLOGY C1 session (row 29) from the
#    "Pending" (yellow) look to the "Not Recorded" (pink) look used by the
#    legend at K21:M21, and update its status text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- widen column I ("Students") from 10 to 14 characters -----------------
# Excel's ColumnWidth setter pads by 5/6 of a character internally, so back
# that padding out to land on an exact stored width of 14.
$ws.Columns.Item(9).ColumnWidth = 14 - 5/6

# --- row 2: reorder "Recorded By" list -------------------------------------
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, System"

# --- row 3: reorder "Recorded By" list -------------------------------------
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, System"

# --- row 7 / row 8: Missing / Pending session counters ----------------------
$ws.Range("L7").Value = 1
$ws.Range("L8").Value = 21

# --- row 15: reorder "Recorded By" list + refresh group stats --------------
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("P15").Value = 1
$ws.Range("Q15").Value = 21

# --- row 29: PHYSIOLOGY C1 session 2 is now "Not Recorded" -----------------
# Re-style A29:I29 from the yellow "Pending" look to the pink "Not Recorded"
# look. Start from a cell that already carries the black-text formatting
# (fontId 2) used by "Recorded"/"Pending" rows (A2), paste just its format,
# then recolor the fill to the legend's pink ("Not Recorded"/"Red" swatch at
# L21) and re-apply the centered alignment.
$ws.Range("A2").Copy()
$ws.Range("A29:I29").PasteSpecial(-4122)
$ws.Range("A29:I29").Interior.Color = 12695295
$ws.Range("A29:I29").HorizontalAlignment = -4108
$ws.Range("A29:I29").VerticalAlignment = -4108
$ws.Range("I29").Value = "Not Recorded"

$excel.CutCopyMode = $false
